$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 606.9400000000001  # H15
$ws.Cells.Item(15, 9).Value = 606.9400000000001  # I15
$ws.Cells.Item(15, 11).Value = 1820.82  # K15
$ws.Cells.Item(15, 13).Value = -1651.82  # M15

$ws.Cells.Item(112, 8).Value = 2716.4167  # H112
$ws.Cells.Item(112, 9).Value = 925.25  # I112
$ws.Cells.Item(112, 10).Value = 2940.3125  # J112
$ws.Cells.Item(112, 11).Value = 2775.75  # K112
$ws.Cells.Item(112, 12).Value = 8820.9375  # L112
$ws.Cells.Item(112, 13).Value = -1667.75  # M112
$ws.Cells.Item(112, 14).Value = -11036.9375  # N112

$ws.Cells.Item(127, 8).Value = 988.84  # H127
$ws.Cells.Item(127, 9).Value = 829.3333  # I127
$ws.Cells.Item(127, 10).Value = 999.0213  # J127
$ws.Cells.Item(127, 11).Value = 2487.9999  # K127
$ws.Cells.Item(127, 12).Value = 2997.0639  # L127
$ws.Cells.Item(127, 13).Value = 2472.0001  # M127
$ws.Cells.Item(127, 14).Value = -12917.0639  # N127

$ws.Cells.Item(138, 8).Value = 2635946.5  # H138
$ws.Cells.Item(138, 9).Value = 5265433.5  # I138
$ws.Cells.Item(138, 10).Value = 6459.3423  # J138
$ws.Cells.Item(138, 11).Value = 15796300.5  # K138
$ws.Cells.Item(138, 12).Value = 19378.0269  # L138
$ws.Cells.Item(138, 13).Value = -15791160.5  # M138
$ws.Cells.Item(138, 14).Value = -29658.0269  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1538  # H2
$ws.Cells.Item(2, 9).Value = 1457  # I2
$ws.Cells.Item(2, 10).Value = 1700  # J2
$ws.Cells.Item(2, 11).Value = 1457  # K2
$ws.Cells.Item(2, 12).Value = 1700  # L2
$ws.Cells.Item(2, 13).Value = -1344  # M2
$ws.Cells.Item(2, 14).Value = -1926  # N2

$ws.Cells.Item(32, 8).Value = 18600.734  # H32
$ws.Cells.Item(32, 9).Value = 17644.564  # I32
$ws.Cells.Item(32, 10).Value = 22329.8  # J32
$ws.Cells.Item(32, 11).Value = 17644.564  # K32
$ws.Cells.Item(32, 12).Value = 22329.8  # L32
$ws.Cells.Item(32, 13).Value = -17357.564  # M32
$ws.Cells.Item(32, 14).Value = -22903.8  # N32

$ws.Cells.Item(116, 8).Value = 1538  # H116
$ws.Cells.Item(116, 9).Value = 1457  # I116
$ws.Cells.Item(116, 10).Value = 1700  # J116
$ws.Cells.Item(116, 11).Value = 1457  # K116
$ws.Cells.Item(116, 12).Value = 1700  # L116
$ws.Cells.Item(116, 13).Value = 837  # M116
$ws.Cells.Item(116, 14).Value = -6288  # N116

$ws.Cells.Item(118, 8).Value = 48409  # H118
$ws.Cells.Item(118, 10).Value = 48409  # J118
$ws.Cells.Item(118, 12).Value = 48409  # L118
$ws.Cells.Item(118, 14).Value = -51723  # N118

$ws.Cells.Item(132, 8).Value = 528365.5  # H132
$ws.Cells.Item(132, 9).Value = 626602.9  # I132
$ws.Cells.Item(132, 11).Value = 1879808.7  # K132
$ws.Cells.Item(132, 13).Value = -1877278.7  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1538  # H3
$ws.Cells.Item(3, 9).Value = 1457  # I3
$ws.Cells.Item(3, 10).Value = 1700  # J3
$ws.Cells.Item(3, 11).Value = 1457  # K3
$ws.Cells.Item(3, 12).Value = 1700  # L3
$ws.Cells.Item(3, 13).Value = -1343  # M3
$ws.Cells.Item(3, 14).Value = -1928  # N3

$ws.Cells.Item(105, 8).Value = 3488.6667  # H105
$ws.Cells.Item(105, 9).Value = 3586.4  # I105
$ws.Cells.Item(105, 11).Value = 3586.4  # K105
$ws.Cells.Item(105, 13).Value = -1839.4  # M105

$ws.Cells.Item(107, 8).Value = 28525.75  # H107
$ws.Cells.Item(107, 9).Value = 41710.848  # I107
$ws.Cells.Item(107, 11).Value = 41710.848  # K107
$ws.Cells.Item(107, 13).Value = -39790.848  # M107

$ws.Cells.Item(133, 8).Value = 74950  # H133
$ws.Cells.Item(133, 10).Value = 74950  # J133
$ws.Cells.Item(133, 12).Value = 74950  # L133
$ws.Cells.Item(133, 14).Value = -85070  # N133

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3083.5083  # H31
$ws.Cells.Item(31, 9).Value = 1838.7273  # I31
$ws.Cells.Item(31, 10).Value = 4550.5713  # J31
$ws.Cells.Item(31, 11).Value = 1838.7273  # K31
$ws.Cells.Item(31, 12).Value = 4550.5713  # L31
$ws.Cells.Item(31, 14).Value = -5140.5713  # N31
$ws.Cells.Item(31, 13).Value = -1543.7273  # M31

$ws.Cells.Item(34, 8).Value = 3083.5083  # H34
$ws.Cells.Item(34, 9).Value = 1838.7273  # I34
$ws.Cells.Item(34, 10).Value = 4550.5713  # J34
$ws.Cells.Item(34, 11).Value = 1838.7273  # K34
$ws.Cells.Item(34, 12).Value = 4550.5713  # L34
$ws.Cells.Item(34, 14).Value = -4954.5713  # N34
$ws.Cells.Item(34, 13).Value = -1636.7273  # M34

$ws.Cells.Item(58, 8).Value = 1432.0714  # H58
$ws.Cells.Item(58, 9).Value = 1484.2916  # I58
$ws.Cells.Item(58, 10).Value = 1118.75  # J58
$ws.Cells.Item(58, 11).Value = 1484.2916  # K58
$ws.Cells.Item(58, 12).Value = 1118.75  # L58
$ws.Cells.Item(58, 13).Value = -1281.2916  # M58
$ws.Cells.Item(58, 14).Value = -1524.75  # N58

$ws.Cells.Item(107, 8).Value = 536.9545000000001  # H107
$ws.Cells.Item(107, 9).Value = 456.66666  # I107
$ws.Cells.Item(107, 10).Value = 709  # J107
$ws.Cells.Item(107, 11).Value = 456.66666  # K107
$ws.Cells.Item(107, 12).Value = 709  # L107
$ws.Cells.Item(107, 13).Value = 1463.33334  # M107
$ws.Cells.Item(107, 14).Value = -4549  # N107

$ws.Cells.Item(134, 8).Value = 1587.3478  # H134
$ws.Cells.Item(134, 9).Value = 1429.75  # I134
$ws.Cells.Item(134, 11).Value = 4289.25  # K134
$ws.Cells.Item(134, 13).Value = -1754.25  # M134

$ws.Cells.Item(136, 8).Value = 1432.0714  # H136
$ws.Cells.Item(136, 9).Value = 1484.2916  # I136
$ws.Cells.Item(136, 10).Value = 1118.75  # J136
$ws.Cells.Item(136, 11).Value = 4452.8748  # K136
$ws.Cells.Item(136, 12).Value = 3356.25  # L136
$ws.Cells.Item(136, 13).Value = -1902.8748  # M136
$ws.Cells.Item(136, 14).Value = -8456.25  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 84256.836  # H4
$ws.Cells.Item(4, 9).Value = 125097.75  # I4
$ws.Cells.Item(4, 10).Value = 2575  # J4
$ws.Cells.Item(4, 11).Value = 375293.25  # K4
$ws.Cells.Item(4, 12).Value = 7725  # L4
$ws.Cells.Item(4, 13).Value = -375181.25  # M4
$ws.Cells.Item(4, 14).Value = -7949  # N4

$ws.Cells.Item(80, 8).Value = 5220  # H80
$ws.Cells.Item(80, 9).Value = 8700  # I80
$ws.Cells.Item(80, 10).Value = 3954.5454  # J80
$ws.Cells.Item(80, 11).Value = 26100  # K80
$ws.Cells.Item(80, 12).Value = 11863.6362  # L80
$ws.Cells.Item(80, 13).Value = -25164  # M80
$ws.Cells.Item(80, 14).Value = -13735.6362  # N80

$ws.Cells.Item(83, 8).Value = 5220  # H83
$ws.Cells.Item(83, 9).Value = 8700  # I83
$ws.Cells.Item(83, 10).Value = 3954.5454  # J83
$ws.Cells.Item(83, 11).Value = 78300  # K83
$ws.Cells.Item(83, 12).Value = 35590.9086  # L83
$ws.Cells.Item(83, 13).Value = -73620  # M83
$ws.Cells.Item(83, 14).Value = -44950.9086  # N83

$ws.Cells.Item(112, 8).Value = 2761  # H112
$ws.Cells.Item(112, 10).Value = 3308.9473  # J112
$ws.Cells.Item(112, 12).Value = 9926.841899999999  # L112
$ws.Cells.Item(112, 14).Value = -12142.8419  # N112

$ws.Cells.Item(129, 8).Value = 2943142  # H129
$ws.Cells.Item(129, 9).Value = 950  # I129
$ws.Cells.Item(129, 10).Value = 3335434.2  # J129
$ws.Cells.Item(129, 11).Value = 2850  # K129
$ws.Cells.Item(129, 12).Value = 10006302.6  # L129
$ws.Cells.Item(129, 13).Value = 2150  # M129
$ws.Cells.Item(129, 14).Value = -10016302.6  # N129

$ws.Cells.Item(131, 8).Value = 10640745  # H131
$ws.Cells.Item(131, 10).Value = 12049595  # J131
$ws.Cells.Item(131, 12).Value = 36148785  # L131
$ws.Cells.Item(131, 14).Value = -36158865  # N131

$ws.Cells.Item(140, 8).Value = 2563.25  # H140
$ws.Cells.Item(140, 9).Value = 1124  # I140
$ws.Cells.Item(140, 10).Value = 8032.4  # J140
$ws.Cells.Item(140, 11).Value = 3372  # K140
$ws.Cells.Item(140, 12).Value = 24097.2  # L140
$ws.Cells.Item(140, 13).Value = 1808  # M140
$ws.Cells.Item(140, 14).Value = -34457.2  # N140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 9233769  # H3
$ws.Cells.Item(3, 9).Value = 910818.2  # I3
$ws.Cells.Item(3, 11).Value = 910818.2  # K3
$ws.Cells.Item(3, 13).Value = -910702.2  # M3

$ws.Cells.Item(5, 8).Value = 43004  # H5
$ws.Cells.Item(5, 9).Value = 43004  # I5
$ws.Cells.Item(5, 10).Value = 0  # J5
$ws.Cells.Item(5, 11).Value = 43004  # K5
$ws.Cells.Item(5, 12).Value = 0  # L5
$ws.Cells.Item(5, 13).Value = -42892  # M5
$ws.Cells.Item(5, 14).ClearContents()  # N5

$ws.Cells.Item(36, 8).Value = 5875.4375  # H36
$ws.Cells.Item(36, 9).Value = 2697  # I36
$ws.Cells.Item(36, 10).Value = 6934.9165  # J36
$ws.Cells.Item(36, 11).Value = 2697  # K36
$ws.Cells.Item(36, 12).Value = 6934.9165  # L36
$ws.Cells.Item(36, 13).Value = -2212  # M36
$ws.Cells.Item(36, 14).Value = -7904.9165  # N36

$ws.Cells.Item(42, 8).Value = 18269.23  # H42
$ws.Cells.Item(42, 10).Value = 18269.23  # J42
$ws.Cells.Item(42, 12).Value = 18269.23  # L42
$ws.Cells.Item(42, 14).Value = -19239.23  # N42

$ws.Cells.Item(51, 8).Value = 30230.691  # H51
$ws.Cells.Item(51, 10).Value = 30230.691  # J51
$ws.Cells.Item(51, 12).Value = 30230.691  # L51
$ws.Cells.Item(51, 14).Value = -31248.691  # N51

$ws.Cells.Item(54, 8).Value = 0  # H54
$ws.Cells.Item(54, 10).Value = 0  # J54
$ws.Cells.Item(54, 12).Value = 0  # L54
$ws.Cells.Item(54, 14).ClearContents()  # N54

$ws.Cells.Item(104, 8).Value = 37500  # H104
$ws.Cells.Item(104, 10).Value = 37500  # J104
$ws.Cells.Item(104, 12).Value = 37500  # L104
$ws.Cells.Item(104, 14).Value = -44488  # N104

$ws.Cells.Item(115, 8).Value = 18269.23  # H115
$ws.Cells.Item(115, 10).Value = 18269.23  # J115
$ws.Cells.Item(115, 12).Value = 18269.23  # L115
$ws.Cells.Item(115, 14).Value = -20619.23  # N115

$ws.Cells.Item(122, 8).Value = 4816.8438  # H122
$ws.Cells.Item(122, 9).Value = 4516  # I122
$ws.Cells.Item(122, 10).Value = 5478.7  # J122
$ws.Cells.Item(122, 11).Value = 13548  # K122
$ws.Cells.Item(122, 12).Value = 16436.1  # L122
$ws.Cells.Item(122, 13).Value = -11098  # M122
$ws.Cells.Item(122, 14).Value = -21336.1  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(62, 8).Value = 39800  # H62
$ws.Cells.Item(62, 10).Value = 39800  # J62
$ws.Cells.Item(62, 12).Value = 39800  # L62
$ws.Cells.Item(62, 14).Value = -41048  # N62

$ws.Cells.Item(63, 8).Value = 43733.332  # H63
$ws.Cells.Item(63, 10).Value = 43733.332  # J63
$ws.Cells.Item(63, 12).Value = 43733.332  # L63
$ws.Cells.Item(63, 14).Value = -45231.332  # N63

$ws.Cells.Item(65, 8).Value = 39800  # H65
$ws.Cells.Item(65, 10).Value = 39800  # J65
$ws.Cells.Item(65, 12).Value = 119400  # L65
$ws.Cells.Item(65, 14).Value = -125640  # N65

$ws.Cells.Item(66, 8).Value = 43733.332  # H66
$ws.Cells.Item(66, 10).Value = 43733.332  # J66
$ws.Cells.Item(66, 12).Value = 131199.996  # L66
$ws.Cells.Item(66, 14).Value = -138687.996  # N66

$ws.Cells.Item(123, 8).Value = 58444.777  # H123
$ws.Cells.Item(123, 10).Value = 58444.777  # J123
$ws.Cells.Item(123, 12).Value = 58444.777  # L123
$ws.Cells.Item(123, 14).Value = -68244.777  # N123

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 501500  # H2
$ws.Cells.Item(2, 10).Value = 1000000  # J2
$ws.Cells.Item(2, 12).Value = 1000000  # L2
$ws.Cells.Item(2, 14).Value = -1000224  # N2

$ws.Cells.Item(27, 8).Value = 0  # H27
$ws.Cells.Item(27, 10).Value = 0  # J27
$ws.Cells.Item(27, 12).Value = 0  # L27
$ws.Cells.Item(27, 14).ClearContents()  # N27

$ws.Cells.Item(136, 8).Value = 1838.25  # H136
$ws.Cells.Item(136, 9).Value = 1892.9678  # I136
$ws.Cells.Item(136, 10).Value = 1499  # J136
$ws.Cells.Item(136, 11).Value = 5678.903399999999  # K136
$ws.Cells.Item(136, 12).Value = 4497  # L136
$ws.Cells.Item(136, 13).Value = -3128.903399999999  # M136
$ws.Cells.Item(136, 14).Value = -9597  # N136
